# Update the "File Name" cells (B1:B6) with the file name + function name
# that the file actually documents, wrap the text, size row/column to fit,
# and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "File Name:a.c`nDocument fun name: func1"
$ws.Range("B2").Value = "File Name:a.c`nDocument fun name: func11"
$ws.Range("B3").Value = "File Name:a.c`nDocument fun name: func12"
$ws.Range("B4").Value = "File Name:b.c`nDocument fun name: funcb"
$ws.Range("B5").Value = "File Name:c.c`nDocument fun name: funcc"
$ws.Range("B6").Value = "File Name:c.c`nDocument fun name: funccc"

# Wrap the multi-line text inside the cells.
$ws.Range("B1:B6").WrapText = $true

# Let rows size themselves first, then pin row 1 to its taller height.
$ws.Rows("1:6").AutoFit() | Out-Null
$ws.Rows(1).RowHeight = 55.5

# Widen column B so the wrapped text is readable.
$ws.Columns("B").ColumnWidth = 57.71

# Move the selection like the author left it.
$ws.Range("D3").Select() | Out-Null
